$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7 and 14-16 (columns A=Name, B=Position, C=Team)
$data = @{
    2  = @("Darius Garland", "PG", "Cleveland Cavaliers")
    3  = @("Scotty Pippen Jr.", "PG,SG", "Memphis Grizzlies")
    4  = @("Austin Reaves", "PG,SG", "Los Angeles Lakers")
    5  = @("Stephen Curry", "PG,SG", "Golden State Warriors")
    6  = @("Trey Murphy III", "SF,PF", "New Orleans Pelicans")
    7  = @("Zach Collins", "PF,C", "San Antonio Spurs")
    14 = @("Jalen Johnson", "PF", "Atlanta Hawks")
    15 = @("Tyrese Haliburton", "PG,SG", "Indiana Pacers")
    16 = @("Keegan Murray", "SF,PF", "Sacramento Kings")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
